# Updated cryptos list on Mon Jun 19 03:37:35 UTC 2023 with GitHub Actions
# Applies the refreshed price / volume(1h) snapshot to the existing rows,
# and reflects that Elrond moved above EnergySwap in the ranking (rows 49/50 swap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.413.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.724.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.14%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "

# Row 6
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4896"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.72%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2610"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06203"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.726.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.05%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07011"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.59%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.578"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.19%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5999"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.82%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "

# Row 16
$ws.Range("E16").Value = "  +0.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.418.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "

# Row 18
$ws.Range("E18").Value = "  +0.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007134"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.78%  "

# Row 20
$ws.Range("E20").Value = "  -1.72%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.944.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.66%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.469"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.596"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.22%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.157"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.40%  "

# Row 26
$ws.Range("E26").Value = "  -0.78%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.392"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.29%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.700"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.956"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.14%  "

# Row 31
$ws.Range("E31").Value = "  -1.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.682"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04535"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9996"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.604"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.28%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9952"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6270"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9087"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.27%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.956"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.20%  "

# Row 40
$ws.Range("E40").Value = "  +0.98%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01481"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.73%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.29%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.439"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.59%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3840"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.66%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.706"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.89%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1157"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05366"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.17%  "

# Rows 49/50: ranking swap - Elrond now ranks above EnergySwap
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.14%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.697"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.68%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.241"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.08%  "
